# Set cell H4 on the active worksheet to "ALL" (reuses existing shared string)
# and move/update the active selection from G4 to H4, matching the authored edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = "ALL"
$ws.Range("H4").Select()
